$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9: H,I,K,M
$ws.Range("H9").Value = 61.6
$ws.Range("I9").Value = 61.6
$ws.Range("K9").Value = 61.6
$ws.Range("M9").Value = 107.4

# Row 32: H,I,J,K,L,M,N
$ws.Range("H32").Value = 798
$ws.Range("I32").Value = 797.3333
$ws.Range("J32").Value = 800
$ws.Range("K32").Value = 797.3333
$ws.Range("L32").Value = 800
$ws.Range("M32").Value = -471.3333
$ws.Range("N32").Value = -1452

# Row 58: H,J,L,N
$ws.Range("H58").Value = 1096.9445
$ws.Range("J58").Value = 1253.9231
$ws.Range("L58").Value = 3761.7693
$ws.Range("N58").Value = -4061.7693

# Row 74: H,I,K,M
$ws.Range("H74").Value = 13000
$ws.Range("I74").Value = 13000
$ws.Range("K74").Value = 13000
$ws.Range("M74").Value = -12064

# Row 77: H,I,K,M
$ws.Range("H77").Value = 13000
$ws.Range("I77").Value = 13000
$ws.Range("K77").Value = 65000
$ws.Range("M77").Value = -60320

# Row 88: H,I,J,K,L,M,N
$ws.Range("H88").Value = 916
$ws.Range("I88").Value = 1271.3334
$ws.Range("J88").Value = 649.5
$ws.Range("K88").Value = 1271.3334
$ws.Range("L88").Value = 649.5
$ws.Range("M88").Value = -865.3334
$ws.Range("N88").Value = -1461.5

# Row 91: H,I,J,K,L,M,N
$ws.Range("H91").Value = 916
$ws.Range("I91").Value = 1271.3334
$ws.Range("J91").Value = 649.5
$ws.Range("K91").Value = 1271.3334
$ws.Range("L91").Value = 649.5
$ws.Range("M91").Value = 132.6666
$ws.Range("N91").Value = -3457.5

# Row 95: H,J,L,N
$ws.Range("H95").Value = 18980
$ws.Range("J95").Value = 18980
$ws.Range("L95").Value = 18980
$ws.Range("N95").Value = -24472

# Row 97: H,I,J,K,L,M,N
$ws.Range("H97").Value = 599.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 599.5
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 1798.5
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = -2790.5

# Row 138: H,I,J,K,L,M,N
$ws.Range("H138").Value = 2187
$ws.Range("I138").Value = 1332.6666
$ws.Range("J138").Value = 4750
$ws.Range("K138").Value = 3997.9998
$ws.Range("L138").Value = 14250
$ws.Range("M138").Value = 1142.0002
$ws.Range("N138").Value = -24530

$ws = $wb.Worksheets.Item("ARM")
# Row 61: H,I,K,M
$ws.Range("H61").Value = 2098.6428
$ws.Range("I61").Value = 1580.1818
$ws.Range("K61").Value = 1580.1818
$ws.Range("M61").Value = -1368.1818

# Row 74: H,I,K,M
$ws.Range("H74").Value = 2977.3076
$ws.Range("I74").Value = 2280.0908
$ws.Range("K74").Value = 2280.0908
$ws.Range("M74").Value = -1406.0908

# Row 77: H,I,K,M
$ws.Range("H77").Value = 2977.3076
$ws.Range("I77").Value = 2280.0908
$ws.Range("K77").Value = 11400.454
$ws.Range("M77").Value = -7032.454

# Row 88: H,J,L,N
$ws.Range("H88").Value = 1506
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = ""

# Row 91: H,J,L,N
$ws.Range("H91").Value = 1506
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = ""

# Row 97: H,I,J,K,L,M,N
$ws.Range("H97").Value = 1080.8
$ws.Range("I97").Value = 1161.25
$ws.Range("J97").Value = 759
$ws.Range("K97").Value = 1161.25
$ws.Range("L97").Value = 759
$ws.Range("M97").Value = -665.25
$ws.Range("N97").Value = -1751

# Row 136: H,I,K,M
$ws.Range("H136").Value = 2098.6428
$ws.Range("I136").Value = 1580.1818
$ws.Range("K136").Value = 4740.5454
$ws.Range("M136").Value = -2190.5454

$ws = $wb.Worksheets.Item("BSM")
# Row 82: H,I,K,M
$ws.Range("H82").Value = 30575.125
$ws.Range("I82").Value = 8829.4
$ws.Range("K82").Value = 8829.4
$ws.Range("M82").Value = -8446.4

# Row 85: H,I,K,M
$ws.Range("H85").Value = 30575.125
$ws.Range("I85").Value = 8829.4
$ws.Range("K85").Value = 8829.4
$ws.Range("M85").Value = -7503.4

# Row 88: H,J,L,N
$ws.Range("H88").Value = 23474.75
$ws.Range("J88").Value = 23474.75
$ws.Range("L88").Value = 23474.75
$ws.Range("N88").Value = -24286.75

# Row 91: H,J,L,N
$ws.Range("H91").Value = 23474.75
$ws.Range("J91").Value = 23474.75
$ws.Range("L91").Value = 23474.75
$ws.Range("N91").Value = -26282.75

# Row 94: H,I,K,M
$ws.Range("H94").Value = 1227.5714
$ws.Range("I94").Value = 1227.5714
$ws.Range("K94").Value = 1227.5714
$ws.Range("M94").Value = -776.5714

$ws = $wb.Worksheets.Item("CRP")
# Row 16: H,I,K,M
$ws.Range("H16").Value = 2294.5
$ws.Range("I16").Value = 2294.5
$ws.Range("K16").Value = 2294.5
$ws.Range("M16").Value = -2007.5

# Row 113: H,I,K,M
$ws.Range("H113").Value = 2294.5
$ws.Range("I113").Value = 2294.5
$ws.Range("K113").Value = 2294.5
$ws.Range("M113").Value = -124.5

$ws = $wb.Worksheets.Item("CUL")
# Row 17: H,I,J,K,L,M,N
$ws.Range("H17").Value = 40
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 120
$ws.Range("L17").Value = 120
$ws.Range("M17").Value = 49
$ws.Range("N17").Value = -458

# Row 46: H,J,L,N
$ws.Range("H46").Value = 6725
$ws.Range("J46").Value = 12500
$ws.Range("L46").Value = 37500
$ws.Range("N46").Value = -37682

# Row 60: H,I,J,K,L,M,N
$ws.Range("H60").Value = 1081.8823
$ws.Range("I60").Value = 176.54546
$ws.Range("J60").Value = 2741.6667
$ws.Range("K60").Value = 529.6363799999999
$ws.Range("L60").Value = 8225.000100000001
$ws.Range("M60").Value = -278.6363799999999
$ws.Range("N60").Value = -8727.000100000001

$ws = $wb.Worksheets.Item("GSM")
# Row 5: H,I,K,M
$ws.Range("H5").Value = 212
$ws.Range("I5").Value = 212
$ws.Range("K5").Value = 212
$ws.Range("M5").Value = -100

# Row 57: H,J,L,N
$ws.Range("H57").Value = 14597.667
$ws.Range("J57").Value = 14597.667
$ws.Range("L57").Value = 14597.667
$ws.Range("N57").Value = -16237.667

# Row 80: H,I,J,K,L,M,N
$ws.Range("H80").Value = 1066
$ws.Range("I80").Value = 999.5
$ws.Range("J80").Value = 1099.25
$ws.Range("K80").Value = 999.5
$ws.Range("L80").Value = 1099.25
$ws.Range("M80").Value = -1.5
$ws.Range("N80").Value = -3095.25

# Row 83: H,I,J,K,L,M,N
$ws.Range("H83").Value = 1066
$ws.Range("I83").Value = 999.5
$ws.Range("J83").Value = 1099.25
$ws.Range("K83").Value = 4997.5
$ws.Range("L83").Value = 5496.25
$ws.Range("M83").Value = -5.5
$ws.Range("N83").Value = -15480.25

# Row 97: H,J,L,N
$ws.Range("H97").Value = 1203.6666
$ws.Range("J97").Value = 1455.5
$ws.Range("L97").Value = 1455.5
$ws.Range("N97").Value = -2447.5

# Row 113: H,I,J,K,L,M,N
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = ""

# Row 132: H,I,J,K,L,M,N
$ws.Range("H132").Value = 253006
$ws.Range("I132").Value = 335674.66
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 1007023.98
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1004493.98
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
# Row 16: H,I,K,M
$ws.Range("H16").Value = 2037.3334
$ws.Range("I16").Value = 1445
$ws.Range("K16").Value = 1445
$ws.Range("M16").Value = -1275

# Row 68: H,I,K,M
$ws.Range("H68").Value = 2443.625
$ws.Range("I68").Value = 1258.3334
$ws.Range("K68").Value = 1258.3334
$ws.Range("M68").Value = -509.3334

# Row 71: H,I,K,M
$ws.Range("H71").Value = 2443.625
$ws.Range("I71").Value = 1258.3334
$ws.Range("K71").Value = 6291.666999999999
$ws.Range("M71").Value = -2547.666999999999

# Row 82: H,J,L,N
$ws.Range("H82").Value = 4251.9
$ws.Range("J82").Value = 5152.375
$ws.Range("L82").Value = 5152.375
$ws.Range("N82").Value = -5874.375

# Row 85: H,J,L,N
$ws.Range("H85").Value = 4251.9
$ws.Range("J85").Value = 5152.375
$ws.Range("L85").Value = 5152.375
$ws.Range("N85").Value = -7648.375

# Row 93: H,J,L,N
$ws.Range("H93").Value = 1133.55
$ws.Range("J93").Value = 1099.5
$ws.Range("L93").Value = 1099.5
$ws.Range("N93").Value = -3595.5

# Row 100: H
$ws.Range("H100").Value = 6857

# Row 122: H,I,K,M
$ws.Range("H122").Value = 1477
$ws.Range("I122").Value = 1477
$ws.Range("K122").Value = 4431
$ws.Range("M122").Value = -1981

# Row 132: H,I,J,K,L,M,N
$ws.Range("H132").Value = 3588.875
$ws.Range("I132").Value = 3578.6155
$ws.Range("J132").Value = 3633.3333
$ws.Range("K132").Value = 10735.8465
$ws.Range("L132").Value = 10899.9999
$ws.Range("M132").Value = -8205.8465
$ws.Range("N132").Value = -15959.9999

# Row 134: H,J,L,N
$ws.Range("H134").Value = 24429
$ws.Range("J134").Value = 24429
$ws.Range("L134").Value = 24429
$ws.Range("N134").Value = -34569

# Row 136: H,I,K,M
$ws.Range("H136").Value = 2407.2856
$ws.Range("I136").Value = 2407.2856
$ws.Range("K136").Value = 7221.8568
$ws.Range("M136").Value = -4671.8568

$ws = $wb.Worksheets.Item("WVR")
# Row 15: H,J,L,N
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""

# Row 68: H,J,L,N
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""

# Row 71: H,J,L,N
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""

# Row 107: H,I,J,K,L,M,N
$ws.Range("H107").Value = 1074.4
$ws.Range("I107").Value = 897
$ws.Range("J107").Value = 1192.6666
$ws.Range("K107").Value = 2691
$ws.Range("L107").Value = 3577.9998
$ws.Range("M107").Value = -771
$ws.Range("N107").Value = -7417.9998

# Row 133: H,J,L,N
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120
